# Applies updated "dSF" (column F) values to Sheet1, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new dSF (column F) value
$updates = @{
    2  = -1
    3  = 2
    4  = 3
    5  = -2
    7  = -3
    8  = 6
    9  = -3
    11 = 5
    12 = -3
    13 = 2
    17 = 0
    18 = 4
    19 = 0
    20 = -1
    21 = -3
    22 = -1
    23 = -3
    25 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
